{"js": "// Update the answer cells of the \"three-digit number divided by one-digit\n// number\" worksheet table: each old \"A\u00f7B=C, D\" expression is replaced by a\n// new one. Every old string is unique in the document, so a straight\n// search-and-replace per pair is sufficient and order independent.\nconst replacements = [\n  [\"635\u00f74=158, 3\", \"454\u00f72=227, 0\"],\n  [\"437\u00f74=109, 1\", \"416\u00f73=138, 2\"],\n  [\"884\u00f75=176, 4\", \"560\u00f73=186, 2\"],\n  [\"594\u00f77=84, 6\", \"971\u00f73=323, 2\"],\n  [\"719\u00f75=143, 4\", \"167\u00f72=83, 1\"],\n  [\"562\u00f73=187, 1\", \"747\u00f73=249, 0\"],\n  [\"968\u00f78=121, 0\", \"371\u00f74=92, 3\"],\n  [\"266\u00f73=88, 2\", \"520\u00f77=74, 2\"],\n  [\"412\u00f77=58, 6\", \"603\u00f79=67, 0\"],\n  [\"268\u00f74=67, 0\", \"558\u00f78=69, 6\"],\n  [\"485\u00f75=97, 0\", \"152\u00f77=21, 5\"],\n  [\"412\u00f73=137, 1\", \"878\u00f76=146, 2\"],\n  [\"230\u00f77=32, 6\", \"607\u00f74=151, 3\"],\n  [\"334\u00f79=37, 1\", \"985\u00f78=123, 1\"],\n  [\"205\u00f73=68, 1\", \"286\u00f74=71, 2\"],\n  [\"170\u00f75=34, 0\", \"733\u00f73=244, 1\"],\n  [\"756\u00f78=94, 4\", \"751\u00f72=375, 1\"],\n  [\"975\u00f77=139, 2\", \"791\u00f78=98, 7\"],\n  [\"144\u00f73=48, 0\", \"702\u00f76=117, 0\"],\n  [\"167\u00f76=27, 5\", \"908\u00f79=100, 8\"],\n  [\"925\u00f74=231, 1\", \"480\u00f75=96, 0\"],\n  [\"107\u00f76=17, 5\", \"654\u00f79=72, 6\"],\n  [\"856\u00f72=428, 0\", \"309\u00f74=77, 1\"],\n  [\"373\u00f73=124, 1\", \"174\u00f72=87, 0\"],\n  [\"867\u00f73=289, 0\", \"598\u00f72=299, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answer cells of the \"three-digit number divided by one-digit\n# number\" worksheet table: each old \"A\u00f7B=C, D\" expression is replaced by a\n# new one. Every old string is unique in the document, so a plain\n# Find/Replace per pair is sufficient and order independent.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"635\u00f74=158, 3\", \"454\u00f72=227, 0\"),\n    @(\"437\u00f74=109, 1\", \"416\u00f73=138, 2\"),\n    @(\"884\u00f75=176, 4\", \"560\u00f73=186, 2\"),\n    @(\"594\u00f77=84, 6\", \"971\u00f73=323, 2\"),\n    @(\"719\u00f75=143, 4\", \"167\u00f72=83, 1\"),\n    @(\"562\u00f73=187, 1\", \"747\u00f73=249, 0\"),\n    @(\"968\u00f78=121, 0\", \"371\u00f74=92, 3\"),\n    @(\"266\u00f73=88, 2\", \"520\u00f77=74, 2\"),\n    @(\"412\u00f77=58, 6\", \"603\u00f79=67, 0\"),\n    @(\"268\u00f74=67, 0\", \"558\u00f78=69, 6\"),\n    @(\"485\u00f75=97, 0\", \"152\u00f77=21, 5\"),\n    @(\"412\u00f73=137, 1\", \"878\u00f76=146, 2\"),\n    @(\"230\u00f77=32, 6\", \"607\u00f74=151, 3\"),\n    @(\"334\u00f79=37, 1\", \"985\u00f78=123, 1\"),\n    @(\"205\u00f73=68, 1\", \"286\u00f74=71, 2\"),\n    @(\"170\u00f75=34, 0\", \"733\u00f73=244, 1\"),\n    @(\"756\u00f78=94, 4\", \"751\u00f72=375, 1\"),\n    @(\"975\u00f77=139, 2\", \"791\u00f78=98, 7\"),\n    @(\"144\u00f73=48, 0\", \"702\u00f76=117, 0\"),\n    @(\"167\u00f76=27, 5\", \"908\u00f79=100, 8\"),\n    @(\"925\u00f74=231, 1\", \"480\u00f75=96, 0\"),\n    @(\"107\u00f76=17, 5\", \"654\u00f79=72, 6\"),\n    @(\"856\u00f72=428, 0\", \"309\u00f74=77, 1\"),\n    @(\"373\u00f73=124, 1\", \"174\u00f72=87, 0\"),\n    @(\"867\u00f73=289, 0\", \"598\u00f72=299, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
